$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: mark column C ("Yes") like the other rows in that block ---
$ws.Range("C15").Value = "Yes"

# --- New shared strings must be created in this exact order so the
#     shared-string table indices line up with the target workbook
#     (DirectX Implementation was typed first, then OpenGL Implementation,
#     then the reusable leaf labels, in the order they first appear). ---
$ws.Range("A65").Value = "DirectX Implementation"
$ws.Range("A57").Value = "OpenGL Implementation"
$ws.Range("B59").Value = "GUI"
$ws.Range("B60").Value = "VertexBuffer"
$ws.Range("B61").Value = "IndexBuffer"
$ws.Range("B63").Value = "Shader"
$ws.Range("B62").Value = "VertexArray"
$ws.Range("B64").Value = "Renderer3D"
$ws.Range("B68").Value = "Device"

# --- OpenGL Implementation block (rows 57-64) ---
$ws.Range("B57").Value = "Window"
$ws.Range("C57").Value = "Yes"

$ws.Range("B58").Value = "Input"
$ws.Range("C58").Value = "Yes"

$ws.Range("C59").Value = "Yes"

$ws.Range("C60").Value = "Yes"

$ws.Range("C61").Value = "Yes"

$ws.Range("C62").Value = "Yes"

$ws.Range("C63").Value = "Yes"

$ws.Range("C64").Value = "Yes"

# --- DirectX Implementation block (rows 65-72) ---
$ws.Range("B65").Value = "Window"
$ws.Range("C65").Value = "Yes"

$ws.Range("B66").Value = "Input"
$ws.Range("C66").Value = "Yes"

$ws.Range("B67").Value = "GUI"
$ws.Range("C67").Value = "Yes"

$ws.Range("C68").Value = "Yes"

$ws.Range("B69").Value = "VertexBuffer"
$ws.Range("C69").Value = "Yes"

$ws.Range("B70").Value = "IndexBuffer"

$ws.Range("B71").Value = "Shader"
$ws.Range("C71").Value = "Yes"

$ws.Range("B72").Value = "Renderer3D"

# --- Formatting: copy the centered/merged-header style (style index 3,
#     used by A26, A37, A41, A45, A51, A53) onto the two new header columns ---
$ws.Range("A26").Copy()
$ws.Range("A57:A72").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Merge the new category header cells, same pattern as the existing blocks ---
$ws.Range("A57:A64").Merge()
$ws.Range("A65:A72").Merge()

# --- View state: keep the selection in sync with where editing ended up ---
$ws.Range("C64").Select()
